$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" column (C) for rows 2-15 from 45170 to 45174
for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45170) {
        $cell.Value2 = 45174
    }
}
